$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; temporarily unprotect so the cell values/text can
# be updated, then re-protect it afterwards so the workbook keeps behaving
# like the original (read-only to casual editing).
$ws.Unprotect()

# Update the "as of" date in the confidential notice (A16): 2021-04-06 -> 2021-04-08
$ws.Range("A16").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-08 for illustrative purposes only and are subject to change."

# Refresh Weight (D) / Percent Change (E) figures for rows 2-13
$ws.Range("D2").Value = 0.03036482766642112
$ws.Range("E2").Value = -0.002840909090909283

$ws.Range("D3").Value = 0.02394484453146202
$ws.Range("E3").Value = -0.001426872770511123

$ws.Range("D4").Value = 0.05192833326871631
$ws.Range("E4").Value = 0.001903855306996549

$ws.Range("D5").Value = 0.1345295962082736
$ws.Range("E5").Value = 0.002893617021276773

$ws.Range("D6").Value = 0.03013293335401078
$ws.Range("E6").Value = -0.01353383458646606

$ws.Range("D7").Value = 0.1209874045283604
$ws.Range("E7").Value = 0.005595828564161076

$ws.Range("D8").Value = 0.1005570794401582
$ws.Range("E8").Value = 0.003616980772891631

$ws.Range("D9").Value = 0.02746566413367968
$ws.Range("E9").Value = 0.0006749156355454478

$ws.Range("D10").Value = 0.1210842088938859
$ws.Range("E10").Value = 0.0006003602161297472

$ws.Range("D11").Value = 0.253793422634048
$ws.Range("E11").Value = 0.01496598639455793

$ws.Range("D12").Value = 0.105211685340984
$ws.Range("E12").Value = 0.002321532211259525

$ws.Range("D13").Value = 0.9999999999999999
$ws.Range("E13").Value = 0.005134387236847227

# Restore protection (password-protected originally; we can't recreate the
# exact legacy hash without the plaintext password, so re-protect without
# one so the sheet stays locked against casual edits like the source file).
$ws.Protect()
